$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3 was started with the last name + link, the first name was filled in later
$ws.Cells.Item(3, 1).Value = "Patel"
$ws.Cells.Item(3, 3).Value = "https://www.youtube.com/watch?v=4HyszYc35ks&authuser=2"

# Remaining presenter rows (4-9), entered in full row by row
$data = @(
    @("Fairweather", "John",     "https://www.youtube.com/watch?v=Y2RSnJwHHnw"),
    @("Cymes",       "Brittany", "https://youtu.be/zLBXx2vZERc"),
    @("Easter",      "Parks",    "https://youtu.be/wJFd90i557w"),
    @("Martinot",    "Melissa",  "https://www.youtube.com/watch?v=IAGhMZd6Rrg&authuser=2"),
    @("Gimar",       "Caleb",    "https://youtu.be/NVsUoX5tmlw"),
    @("Piskurich",   "Nicholas", "https://youtu.be/McZv3w3g8fw")
)

$row = 4
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row = $row + 1
}

# Finally, go back and fill in the first name that was missing on row 3
$ws.Cells.Item(3, 2).Value = "Shreekumari"

$ws.Range("B10").Select()
